$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 300
$ws1.Range("F4").Value = 1184
$ws1.Range("F5").Value = 601

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 300
$ws4.Range("F4").Value = 1184
$ws4.Range("F6").Value = 601
